$d = $word.ActiveDocument

# The paragraph currently ends in: ... <<testprop2>> is another.
# We need to turn that into:       ... <<testprop2>> is another, while <<testprop1>> is a third.
#
# Step 1: strip the trailing period from " is another." -> " is another"
$d.Content.Find.Execute(" is another.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " is another", 2)

# Step 2: append a new run ", while " right after the text we just edited.
$end = $d.Content.End
$insertionPoint = $d.Range($end, $end)
$insertionPoint.InsertAfter(", while ")

# Step 3: append a new simple field referencing testprop1 (same DOCPROPERTY
# field used earlier in the paragraph), as its own fldSimple element.
$end = $d.Content.End
$insertionPoint = $d.Range($end, $end)
$d.Fields.Add($insertionPoint, 1, " DOCPROPERTY  testprop1  \* MERGEFORMAT ", $false)

# Step 4: append the closing run " is a third."
$end = $d.Content.End
$insertionPoint = $d.Range($end, $end)
$insertionPoint.InsertAfter(" is a third.")
